# Auto-applied recomputed-price update to the Leve profit sheets (H-N columns).
# Values come from a refreshed Universalis price snapshot; row/col layout is unchanged.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 76: H76,I76,J76,K76,L76,M76,N76
$ws.Range("H76").Value = 3850
$ws.Range("I76").Value = 2700
$ws.Range("J76").Value = 4080
$ws.Range("K76").Value = 2700
$ws.Range("L76").Value = 4080
$ws.Range("M76").Value = -2385
$ws.Range("N76").Value = -4710

# Row 79: H79,I79,J79,K79,L79,M79,N79
$ws.Range("H79").Value = 3850
$ws.Range("I79").Value = 2700
$ws.Range("J79").Value = 4080
$ws.Range("K79").Value = 2700
$ws.Range("L79").Value = 4080
$ws.Range("M79").Value = -1608
$ws.Range("N79").Value = -6264

# Row 86: H86,J86,L86,N86
$ws.Range("H86").Value = 4077.923
$ws.Range("J86").Value = 5864.4165
$ws.Range("L86").Value = 5864.4165
$ws.Range("N86").Value = -8110.4165

# Row 89: H89,J89,L89,N89
$ws.Range("H89").Value = 4077.923
$ws.Range("J89").Value = 5864.4165
$ws.Range("L89").Value = 29322.0825
$ws.Range("N89").Value = -40554.0825

# Row 112: H112,J112,L112,N112
$ws.Range("H112").Value = 1034.5294
$ws.Range("J112").Value = 1045.8
$ws.Range("L112").Value = 3137.4
$ws.Range("N112").Value = -5353.4

# Row 129: H129,J129,L129,N129
$ws.Range("H129").Value = 4017.25
$ws.Range("J129").Value = 1020.4783
$ws.Range("L129").Value = 3061.4349
$ws.Range("N129").Value = -13061.4349

$ws = $wb.Worksheets.Item("ARM")
# Row 61: H61,I61,J61,K61,L61,M61,N61
$ws.Range("H61").Value = 2472.7932
$ws.Range("I61").Value = 2000
$ws.Range("J61").Value = 2596.1304
$ws.Range("K61").Value = 2000
$ws.Range("L61").Value = 2596.1304
$ws.Range("M61").Value = -1788
$ws.Range("N61").Value = -3020.1304

# Row 74: H74,I74,J74,K74,L74,M74,N74
$ws.Range("H74").Value = 2877.75
$ws.Range("I74").Value = 2440.4
$ws.Range("J74").Value = 3190.1428
$ws.Range("K74").Value = 2440.4
$ws.Range("L74").Value = 3190.1428
$ws.Range("M74").Value = -1566.4
$ws.Range("N74").Value = -4938.1428

# Row 77: H77,I77,J77,K77,L77,M77,N77
$ws.Range("H77").Value = 2877.75
$ws.Range("I77").Value = 2440.4
$ws.Range("J77").Value = 3190.1428
$ws.Range("K77").Value = 12202
$ws.Range("L77").Value = 15950.714
$ws.Range("M77").Value = -7834
$ws.Range("N77").Value = -24686.714

# Row 110: H110,I110,J110,K110,L110,M110,N110
$ws.Range("H110").Value = 27835314
$ws.Range("I110").Value = 35786956
$ws.Range("J110").Value = 4574.75
$ws.Range("K110").Value = 35786956
$ws.Range("L110").Value = 4574.75
$ws.Range("M110").Value = -35784911
$ws.Range("N110").Value = -8664.75

# Row 122: H122,I122,K122,M122
$ws.Range("H122").Value = 1455.8
$ws.Range("I122").Value = 1334.0667
$ws.Range("K122").Value = 4002.2001
$ws.Range("M122").Value = -1552.2001

# Row 132: H132,I132,J132,K132,L132,M132,N132
$ws.Range("H132").Value = 3277.7144
$ws.Range("I132").Value = 3724.85
$ws.Range("J132").Value = 1290.4445
$ws.Range("K132").Value = 11174.55
$ws.Range("L132").Value = 3871.3335
$ws.Range("M132").Value = -8644.549999999999
$ws.Range("N132").Value = -8931.333500000001

# Row 136: H136,I136,J136,K136,L136,M136,N136
$ws.Range("H136").Value = 2472.7932
$ws.Range("I136").Value = 2000
$ws.Range("J136").Value = 2596.1304
$ws.Range("K136").Value = 6000
$ws.Range("L136").Value = 7788.3912
$ws.Range("M136").Value = -3450
$ws.Range("N136").Value = -12888.3912

$ws = $wb.Worksheets.Item("BSM")
# Row 86: H86,I86,J86,K86,L86,M86,N86
$ws.Range("H86").Value = 57284.95
$ws.Range("I86").Value = 101992
$ws.Range("J86").Value = 2643
$ws.Range("K86").Value = 101992
$ws.Range("L86").Value = 2643
$ws.Range("M86").Value = -100869
$ws.Range("N86").Value = -4889

# Row 89: H89,I89,J89,K89,L89,M89,N89
$ws.Range("H89").Value = 57284.95
$ws.Range("I89").Value = 101992
$ws.Range("J89").Value = 2643
$ws.Range("K89").Value = 509960
$ws.Range("L89").Value = 13215
$ws.Range("M89").Value = -504344
$ws.Range("N89").Value = -24447

# Row 116: H116,J116,L116,N116
$ws.Range("H116").Value = 46323.332
$ws.Range("J116").Value = 46323.332
$ws.Range("L116").Value = 46323.332
$ws.Range("N116").Value = -55501.332

$ws = $wb.Worksheets.Item("CRP")
# Row 31: H31,J31,L31,N31
$ws.Range("H31").Value = 18923.775
$ws.Range("J31").Value = 2000.2972
$ws.Range("L31").Value = 2000.2972
$ws.Range("N31").Value = -2590.2972

# Row 34: H34,J34,L34,N34
$ws.Range("H34").Value = 18923.775
$ws.Range("J34").Value = 2000.2972
$ws.Range("L34").Value = 2000.2972
$ws.Range("N34").Value = -2404.2972

# Row 99: H99,I99,J99,K99,L99,M99,N99
$ws.Range("H99").Value = 8122.2
$ws.Range("I99").Value = 3203.4546
$ws.Range("J99").Value = 14134
$ws.Range("K99").Value = 3203.4546
$ws.Range("L99").Value = 14134
$ws.Range("M99").Value = -1705.4546
$ws.Range("N99").Value = -17130

# Row 122: H122,I122,J122,K122,L122,M122,N122
$ws.Range("H122").Value = 2224.2424
$ws.Range("I122").Value = 2092.3572
$ws.Range("J122").Value = 2962.8
$ws.Range("K122").Value = 6277.071599999999
$ws.Range("L122").Value = 8888.400000000001
$ws.Range("M122").Value = -3827.071599999999
$ws.Range("N122").Value = -13788.4

# Row 126: H126,I126,J126,K126,L126,M126,N126
$ws.Range("H126").Value = 8122.2
$ws.Range("I126").Value = 3203.4546
$ws.Range("J126").Value = 14134
$ws.Range("K126").Value = 9610.363799999999
$ws.Range("L126").Value = 42402
$ws.Range("M126").Value = -7140.363799999999
$ws.Range("N126").Value = -47342

$ws = $wb.Worksheets.Item("CUL")
# Row 34: H34,J34,L34,N34
$ws.Range("H34").Value = 1799.875
$ws.Range("J34").Value = 2349.8333
$ws.Range("L34").Value = 7049.499899999999
$ws.Range("N34").Value = -7217.499899999999

# Row 39: H39,J39,L39,N39
$ws.Range("H39").Value = 60000
$ws.Range("J39").Value = 60000
$ws.Range("L39").Value = 180000
$ws.Range("N39").Value = -180588

# Row 55: H55,J55,L55,N55
$ws.Range("H55").Value = 10725.75
$ws.Range("J55").Value = 11646.272
$ws.Range("L55").Value = 34938.81600000001
$ws.Range("N55").Value = -35292.81600000001

# Row 131: H131,I131,J131,K131,L131,M131,N131
$ws.Range("H131").Value = 837.3200000000001
$ws.Range("I131").Value = 660
$ws.Range("J131").Value = 842.80414
$ws.Range("K131").Value = 1980
$ws.Range("L131").Value = 2528.41242
$ws.Range("M131").Value = 3060
$ws.Range("N131").Value = -12608.41242

$ws = $wb.Worksheets.Item("GSM")
# Row 62: H62,I62,K62,M62
$ws.Range("H62").Value = 15000
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").Value = $null

# Row 65: H65,I65,K65,M65
$ws.Range("H65").Value = 15000
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").Value = $null

# Row 132: H132,I132,J132,K132,L132,M132,N132
$ws.Range("H132").Value = 2795.261
$ws.Range("I132").Value = 2752.2104
$ws.Range("J132").Value = 2999.75
$ws.Range("K132").Value = 8256.6312
$ws.Range("L132").Value = 8999.25
$ws.Range("M132").Value = -5726.6312
$ws.Range("N132").Value = -14059.25

# Row 134: H134,J134,L134,N134
$ws.Range("H134").Value = 19079.4
$ws.Range("J134").Value = 19079.4
$ws.Range("L134").Value = 57238.2
$ws.Range("N134").Value = -62308.2

$ws = $wb.Worksheets.Item("LTW")
# Row 40: H40,I40,K40,M40
$ws.Range("H40").Value = 92821.63
$ws.Range("I40").Value = 112617.555
$ws.Range("K40").Value = 112617.555
$ws.Range("M40").Value = -112481.555

# Row 68: H68,I68,J68,K68,L68,M68,N68
$ws.Range("H68").Value = 4275.7144
$ws.Range("I68").Value = 1250
$ws.Range("J68").Value = 5486
$ws.Range("K68").Value = 1250
$ws.Range("L68").Value = 5486
$ws.Range("M68").Value = -501
$ws.Range("N68").Value = -6984

# Row 71: H71,I71,J71,K71,L71,M71,N71
$ws.Range("H71").Value = 4275.7144
$ws.Range("I71").Value = 1250
$ws.Range("J71").Value = 5486
$ws.Range("K71").Value = 6250
$ws.Range("L71").Value = 27430
$ws.Range("M71").Value = -2506
$ws.Range("N71").Value = -34918

# Row 122: H122,I122,J122,K122,L122,M122,N122
$ws.Range("H122").Value = 3392.5833
$ws.Range("I122").Value = 2958.7144
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 8876.143199999999
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -6426.143199999999
$ws.Range("N122").Value = -16900

# Row 132: H132,I132,J132,K132,L132,M132,N132
$ws.Range("H132").Value = 17817.334
$ws.Range("I132").Value = 17817.334
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 53452.00199999999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -50922.00199999999
$ws.Range("N132").Value = $null

# Row 134: H134,J134,L134,N134
$ws.Range("H134").Value = 60864.875
$ws.Range("J134").Value = 60864.875
$ws.Range("L134").Value = 60864.875
$ws.Range("N134").Value = -71004.875

# Row 135: H135,J135,L135,N135
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").Value = $null

# Row 136: H136,I136,J136,K136,L136,M136,N136
$ws.Range("H136").Value = 1254.2222
$ws.Range("I136").Value = 1081.8334
$ws.Range("J136").Value = 2633.3333
$ws.Range("K136").Value = 3245.5002
$ws.Range("L136").Value = 7899.999899999999
$ws.Range("M136").Value = -695.5001999999999
$ws.Range("N136").Value = -12999.9999

# Row 138: H138,J138,L138,N138
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").Value = $null

# Row 139: H139,J139,L139,N139
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").Value = $null

# Row 141: H141,J141,L141,N141
$ws.Range("H141").Value = 65712.5
$ws.Range("J141").Value = 65712.5
$ws.Range("L141").Value = 65712.5
$ws.Range("N141").Value = -76072.5

$ws = $wb.Worksheets.Item("WVR")
# Row 41: H41,J41,L41,M41
$ws.Range("H41").Value = 42459
$ws.Range("J41").Value = 42459
$ws.Range("L41").Value = 42459
$ws.Range("M41").Value = -43239

# Row 122: H122,I122,K122,M122
$ws.Range("H122").Value = 1020.2857
$ws.Range("I122").Value = 1025.3334
$ws.Range("K122").Value = 3076.0002
$ws.Range("M122").Value = -626.0001999999999

# Row 126: H126,I126,J126,K126,L126,M126,N126
$ws.Range("H126").Value = 1526.4736
$ws.Range("I126").Value = 1652.1538
$ws.Range("J126").Value = 1254.1666
$ws.Range("K126").Value = 4956.4614
$ws.Range("L126").Value = 3762.4998
$ws.Range("M126").Value = -2486.4614
$ws.Range("N126").Value = -8702.4998
